# Fix the typo "wir davor" -> "wie davor" in the "Sonstiges" paragraph and
# move the hidden "_GoBack" bookmark from the start of the following
# paragraph to the point (inside this paragraph) where the user's cursor
# last was, i.e. right before "ingeplant".
#
# The resulting run layout (run splits + bookmark placement) mirrors the
# exact structure recorded in the target OOXML.

$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (currently sitting right before
#    "Außerdem waren die Serverkosten ..." in the next paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Re-create "_GoBack" at its new location, right before "ingeplant", in
#    the "Die letzten 2 Wochen ..." paragraph. Doing this first (before any
#    other edit touches the paragraph) keeps the trailing " " run - that
#    follows "... eingeplant wurden." - from being coalesced back into the
#    preceding text run.
$rGoBack = $d.Content
$rGoBack.Find.Execute("ingeplant", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rGoBack.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rGoBack)

# 3) Split off "Die letzten 2 Wochen w" from the rest of the paragraph with
#    a transient bookmark - the split survives even after the bookmark
#    that produced it is removed again.
$rSplitA = $d.Content
$rSplitA.Find.Execute("Die letzten 2 Wochen w", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rSplitA.Collapse(0)
$d.Bookmarks.Add("ZZZ_Split1", $rSplitA)

# 4) Locate the typo "wir" (as part of the unique phrase "wir davor" so we
#    do not match an unrelated "wir" substring elsewhere, e.g. in "wird"),
#    then correct it to "wie".
$rPhrase = $d.Content
$rPhrase.Find.Execute("wir davor", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rTypo = $d.Range($rPhrase.Start, $rPhrase.Start + 3)
$rTypo.Text = "wie"

# 5) Split off the corrected word from what follows it with another
#    transient bookmark.
$rTypo.Collapse(0)
$d.Bookmarks.Add("ZZZ_Split2", $rTypo)

# 6) Drop the transient bookmarks again - only "_GoBack" should remain in
#    the saved document.
$d.Bookmarks.Item("ZZZ_Split1").Delete()
$d.Bookmarks.Item("ZZZ_Split2").Delete()
